$wb = $excel.ActiveWorkbook

# Grab a reference to an existing header cell so we can clone its exact style
# (bold font, thin border, centered alignment) onto the new sheet's header row.
$styleSource = $wb.Worksheets.Item(1).Range("A1")

# Add new worksheet "ODI Batting Extra" after the last existing sheet
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "ODI Batting Extra"

$newSheet.Range("A1").Value = "MATCH_CODE"
$newSheet.Range("B1").Value = "BATTING_POSITION"
$newSheet.Range("C1").Value = "NUM_4"
$newSheet.Range("D1").Value = "NUM_6"
$newSheet.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$newSheet.Range("F1").Value = "MAN_OF_MATCH"

# Copy formatting only (bold + border + centered) from the existing header style
$styleSource.Copy()
$newSheet.Range("A1:F1").PasteSpecial(-4122)

# Row 2 data: force text storage (so "4517" isn't coerced to a number) and then
# drop the temporary Text number-format so the cell keeps the default style.
$newSheet.Range("A2").NumberFormat = "@"
$newSheet.Range("A2").Value = "4517"
$newSheet.Range("A2").ClearFormats()

# These columns are present but blank for this row - still create the cells.
$newSheet.Range("B2:E2").NumberFormat = "@"
$newSheet.Range("B2:E2").Value = ""
$newSheet.Range("B2:E2").ClearFormats()

$newSheet.Range("F2").Value = "NO"

# Restore the originally active sheet/tab.
$wb.Worksheets.Item(1).Activate()
